# New crime data collected - weekly update for cs-en-us-022pct
# Updates header volume/issue number, reporting week dates, and the
# weekly crime-complaint statistics table (Central Park Precinct).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: volume/number and reporting week dates
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# ---------------------------------------------------------------------
# Helper template cells for the "no data" text markers used throughout
# the table: shared string "0" (style 14) and "***.*" (style 14).
# ---------------------------------------------------------------------
$zeroTemplate = $ws.Range("C14")
$naTemplate = $ws.Range("E14")

# ---------------------------------------------------------------------
# Row 15 (Murder)
# ---------------------------------------------------------------------
$zeroTemplate.Copy($ws.Range("D15"))
$naTemplate.Copy($ws.Range("E15"))
$ws.Range("L15").Value = -66.666666666666

# ---------------------------------------------------------------------
# Row 16 (Rape)
# ---------------------------------------------------------------------
$zeroTemplate.Copy($ws.Range("D16"))
$naTemplate.Copy($ws.Range("E16"))
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -80
$ws.Range("N16").Value = -90.055248618784

# ---------------------------------------------------------------------
# Row 18 (Fel. Assault)
# ---------------------------------------------------------------------
$ws.Range("M18").Value = 66.666666666666

# ---------------------------------------------------------------------
# Row 19 (Burglary)
# ---------------------------------------------------------------------
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = -66.666666666666
$ws.Range("J19").Value = 25
$ws.Range("K19").Value = 80
$ws.Range("N19").Value = -72.392638036809

# ---------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------
$ws.Range("D21").Value = 1
$ws.Range("G21").Value = 9
$ws.Range("H21").Value = -77.777777777777
$ws.Range("J21").Value = 72
$ws.Range("K21").Value = 8.333333333333
$ws.Range("L21").Value = 47.169811320754
$ws.Range("M21").Value = -18.75
$ws.Range("N21").Value = -81.986143187067

# ---------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 3
$zeroTemplate.Copy($ws.Range("D24"))
$naTemplate.Copy($ws.Range("E24"))
$ws.Range("F24").Value = 5
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 150
$ws.Range("I24").Value = 38
$ws.Range("K24").Value = 31.034482758620
$ws.Range("L24").Value = 26.666666666666
$ws.Range("M24").Value = -59.139784946236

# ---------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------
$zeroTemplate.Copy($ws.Range("D25"))
$naTemplate.Copy($ws.Range("E25"))
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 100

# ---------------------------------------------------------------------
# Row 26 (UCR Rape*)
# ---------------------------------------------------------------------
$zeroTemplate.Copy($ws.Range("D26"))
$naTemplate.Copy($ws.Range("E26"))
$ws.Range("L26").Value = -75

# ---------------------------------------------------------------------
# Row 27 (Other Sex Crimes)
# ---------------------------------------------------------------------
$zeroTemplate.Copy($ws.Range("C27"))

# ---------------------------------------------------------------------
# Row 30 (Hate Crimes)
# ---------------------------------------------------------------------
$ws.Range("F16").Copy($ws.Range("F30"))
$ws.Range("I30").Value = 3
$ws.Range("L30").Value = 50
